$d = $word.ActiveDocument

# --- Fix spelling: "maximale" -> "maksimale" -------------------------------
# "maximale" is unique in the document, so a plain Find/Replace is safe and
# leaves $r re-seated over the freshly-inserted replacement text.
$r = $d.Content
$found = $r.Find.Execute("maximale", $true, $false, $false, $false, $false, $true, 1, $false, "maksimale", 2)

if (-not $found) {
    throw "Could not find 'maximale' to correct to 'maksimale'."
}

# --- Relocate the "_GoBack" bookmark ---------------------------------------
# Word tracks the most recent edit location with a single "_GoBack"
# bookmark; adding it again moves it (and implicitly drops the old one,
# wherever it used to be - in this document that was at the very end of the
# "... madplan for ugen." paragraph). The edit above landed the cursor right
# after "maks" inside the corrected word, so that's where it goes now.
$splitPos = $r.Start + 4
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
